# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The worker/period detail table occupies B16:J47 (row 15 is the header,
# row 48 is a separately-styled closing row that must stay put).
# The sheet was re-sorted ascending by "Periodo Mora" (column E), keeping
# the relative order of rows that share the same period (stable sort).
$sortRange = $ws.Range("B16:J47")
$keyRange = $ws.Range("E16:E47")

$sortRange.Sort(
    $keyRange,
    1,                              # xlAscending
    [System.Type]::Missing,
    [System.Type]::Missing,
    1,                              # xlAscending (key2, unused)
    [System.Type]::Missing,
    1,                              # xlAscending (key3, unused)
    2,                              # xlTopToBottom
    [System.Type]::Missing,
    1,                              # xlSortNormal
    1,                              # xlTopToBottom (orientation)
    1,                              # xlSortLabels / header=no (range has no header row)
    [System.Type]::Missing,
    [System.Type]::Missing,
    [System.Type]::Missing
)

# Correct the basic salary (Salario Basico) for MARIANA ALEJANDRA FORTICH
# AYOLA's row (now at row 47 after the sort, since her period 1708 is the
# highest) from 1300000 to 737717.
$ws.Range("G47").Value = 737717
